# Applies the cell-value updates captured in the commit's OOXML diff
# (Leve profitability recalculation across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).
# Each entry is Sheet/Row/Col (1-based) plus the new value; Value = $null means the
# cell is cleared entirely (matches the diff removing that <c> element, e.g. when a
# profit computation becomes undefined).

$wb = $excel.ActiveWorkbook

$changes = @(
    @{Sheet="ALC"; Row=19; Col=8; Value=614.75}
    @{Sheet="ALC"; Row=19; Col=9; Value=0.0}
    @{Sheet="ALC"; Row=19; Col=10; Value=614.75}
    @{Sheet="ALC"; Row=19; Col=11; Value=0.0}
    @{Sheet="ALC"; Row=19; Col=12; Value=614.75}
    @{Sheet="ALC"; Row=19; Col=13; Value=$null}
    @{Sheet="ALC"; Row=19; Col=14; Value=-964.75}
    @{Sheet="ALC"; Row=137; Col=8; Value=4210.7104}
    @{Sheet="ALC"; Row=137; Col=9; Value=1061.8235}
    @{Sheet="ALC"; Row=137; Col=10; Value=6759.8096}
    @{Sheet="ALC"; Row=137; Col=11; Value=3185.4705}
    @{Sheet="ALC"; Row=137; Col=12; Value=20279.4288}
    @{Sheet="ALC"; Row=137; Col=13; Value=-635.4704999999999}
    @{Sheet="ALC"; Row=137; Col=14; Value=-25379.4288}
    @{Sheet="ALC"; Row=140; Col=8; Value=39000.0}
    @{Sheet="ALC"; Row=140; Col=10; Value=39000.0}
    @{Sheet="ALC"; Row=140; Col=12; Value=39000.0}
    @{Sheet="ALC"; Row=140; Col=14; Value=-49360.0}
    @{Sheet="ARM"; Row=61; Col=8; Value=31253056.0}
    @{Sheet="ARM"; Row=61; Col=10; Value=3444.25}
    @{Sheet="ARM"; Row=61; Col=12; Value=3444.25}
    @{Sheet="ARM"; Row=61; Col=14; Value=-3868.25}
    @{Sheet="ARM"; Row=74; Col=8; Value=8877.538}
    @{Sheet="ARM"; Row=74; Col=9; Value=1020.75}
    @{Sheet="ARM"; Row=74; Col=10; Value=21448.4}
    @{Sheet="ARM"; Row=74; Col=11; Value=1020.75}
    @{Sheet="ARM"; Row=74; Col=12; Value=21448.4}
    @{Sheet="ARM"; Row=74; Col=13; Value=-146.75}
    @{Sheet="ARM"; Row=74; Col=14; Value=-23196.4}
    @{Sheet="ARM"; Row=77; Col=8; Value=8877.538}
    @{Sheet="ARM"; Row=77; Col=9; Value=1020.75}
    @{Sheet="ARM"; Row=77; Col=10; Value=21448.4}
    @{Sheet="ARM"; Row=77; Col=11; Value=5103.75}
    @{Sheet="ARM"; Row=77; Col=12; Value=107242.0}
    @{Sheet="ARM"; Row=77; Col=13; Value=-735.75}
    @{Sheet="ARM"; Row=77; Col=14; Value=-115978.0}
    @{Sheet="ARM"; Row=132; Col=8; Value=7768.1387}
    @{Sheet="ARM"; Row=132; Col=9; Value=5317.567}
    @{Sheet="ARM"; Row=132; Col=10; Value=20021.0}
    @{Sheet="ARM"; Row=132; Col=11; Value=15952.701}
    @{Sheet="ARM"; Row=132; Col=12; Value=60063.0}
    @{Sheet="ARM"; Row=132; Col=13; Value=-13422.701}
    @{Sheet="ARM"; Row=132; Col=14; Value=-65123.0}
    @{Sheet="ARM"; Row=136; Col=8; Value=31253056.0}
    @{Sheet="ARM"; Row=136; Col=10; Value=3444.25}
    @{Sheet="ARM"; Row=136; Col=12; Value=10332.75}
    @{Sheet="ARM"; Row=136; Col=14; Value=-15432.75}
    @{Sheet="BSM"; Row=80; Col=8; Value=952.0}
    @{Sheet="BSM"; Row=80; Col=9; Value=1447.5}
    @{Sheet="BSM"; Row=80; Col=10; Value=786.8333}
    @{Sheet="BSM"; Row=80; Col=11; Value=1447.5}
    @{Sheet="BSM"; Row=80; Col=12; Value=786.8333}
    @{Sheet="BSM"; Row=80; Col=13; Value=-449.5}
    @{Sheet="BSM"; Row=80; Col=14; Value=-2782.8333}
    @{Sheet="BSM"; Row=83; Col=8; Value=952.0}
    @{Sheet="BSM"; Row=83; Col=9; Value=1447.5}
    @{Sheet="BSM"; Row=83; Col=10; Value=786.8333}
    @{Sheet="BSM"; Row=83; Col=11; Value=7237.5}
    @{Sheet="BSM"; Row=83; Col=12; Value=3934.1665}
    @{Sheet="BSM"; Row=83; Col=13; Value=-2245.5}
    @{Sheet="BSM"; Row=83; Col=14; Value=-13918.1665}
    @{Sheet="BSM"; Row=99; Col=8; Value=0.0}
    @{Sheet="BSM"; Row=99; Col=9; Value=0.0}
    @{Sheet="BSM"; Row=99; Col=10; Value=0.0}
    @{Sheet="BSM"; Row=99; Col=11; Value=0.0}
    @{Sheet="BSM"; Row=99; Col=12; Value=0.0}
    @{Sheet="BSM"; Row=99; Col=13; Value=$null}
    @{Sheet="BSM"; Row=99; Col=14; Value=$null}
    @{Sheet="BSM"; Row=134; Col=8; Value=1560.1818}
    @{Sheet="BSM"; Row=134; Col=9; Value=1476.909}
    @{Sheet="BSM"; Row=134; Col=11; Value=4430.727000000001}
    @{Sheet="BSM"; Row=134; Col=13; Value=-1895.727000000001}
    @{Sheet="CRP"; Row=31; Col=8; Value=23838038.0}
    @{Sheet="CRP"; Row=31; Col=9; Value=50001716.0}
    @{Sheet="CRP"; Row=31; Col=10; Value=52875.91}
    @{Sheet="CRP"; Row=31; Col=11; Value=50001716.0}
    @{Sheet="CRP"; Row=31; Col=12; Value=52875.91}
    @{Sheet="CRP"; Row=31; Col=13; Value=-50001421.0}
    @{Sheet="CRP"; Row=31; Col=14; Value=-53465.91}
    @{Sheet="CRP"; Row=34; Col=8; Value=23838038.0}
    @{Sheet="CRP"; Row=34; Col=9; Value=50001716.0}
    @{Sheet="CRP"; Row=34; Col=10; Value=52875.91}
    @{Sheet="CRP"; Row=34; Col=11; Value=50001716.0}
    @{Sheet="CRP"; Row=34; Col=12; Value=52875.91}
    @{Sheet="CRP"; Row=34; Col=13; Value=-50001514.0}
    @{Sheet="CRP"; Row=34; Col=14; Value=-53279.91}
    @{Sheet="CRP"; Row=58; Col=8; Value=1068.8077}
    @{Sheet="CRP"; Row=58; Col=9; Value=972.27905}
    @{Sheet="CRP"; Row=58; Col=10; Value=1530.0}
    @{Sheet="CRP"; Row=58; Col=11; Value=972.27905}
    @{Sheet="CRP"; Row=58; Col=12; Value=1530.0}
    @{Sheet="CRP"; Row=58; Col=13; Value=-769.27905}
    @{Sheet="CRP"; Row=58; Col=14; Value=-1936.0}
    @{Sheet="CRP"; Row=99; Col=8; Value=5000.0}
    @{Sheet="CRP"; Row=99; Col=9; Value=0.0}
    @{Sheet="CRP"; Row=99; Col=10; Value=5000.0}
    @{Sheet="CRP"; Row=99; Col=11; Value=0.0}
    @{Sheet="CRP"; Row=99; Col=12; Value=5000.0}
    @{Sheet="CRP"; Row=99; Col=13; Value=$null}
    @{Sheet="CRP"; Row=99; Col=14; Value=-7996.0}
    @{Sheet="CRP"; Row=126; Col=8; Value=5000.0}
    @{Sheet="CRP"; Row=126; Col=9; Value=0.0}
    @{Sheet="CRP"; Row=126; Col=10; Value=5000.0}
    @{Sheet="CRP"; Row=126; Col=11; Value=0.0}
    @{Sheet="CRP"; Row=126; Col=12; Value=15000.0}
    @{Sheet="CRP"; Row=126; Col=13; Value=$null}
    @{Sheet="CRP"; Row=126; Col=14; Value=-19940.0}
    @{Sheet="CRP"; Row=132; Col=8; Value=19611812.0}
    @{Sheet="CRP"; Row=132; Col=9; Value=29416494.0}
    @{Sheet="CRP"; Row=132; Col=11; Value=88249482.0}
    @{Sheet="CRP"; Row=132; Col=13; Value=-88246952.0}
    @{Sheet="CRP"; Row=134; Col=8; Value=2036.4333}
    @{Sheet="CRP"; Row=134; Col=9; Value=1810.6}
    @{Sheet="CRP"; Row=134; Col=11; Value=5431.799999999999}
    @{Sheet="CRP"; Row=134; Col=13; Value=-2896.799999999999}
    @{Sheet="CRP"; Row=136; Col=8; Value=1068.8077}
    @{Sheet="CRP"; Row=136; Col=9; Value=972.27905}
    @{Sheet="CRP"; Row=136; Col=10; Value=1530.0}
    @{Sheet="CRP"; Row=136; Col=11; Value=2916.83715}
    @{Sheet="CRP"; Row=136; Col=12; Value=4590.0}
    @{Sheet="CRP"; Row=136; Col=13; Value=-366.8371499999998}
    @{Sheet="CRP"; Row=136; Col=14; Value=-9690.0}
    @{Sheet="CUL"; Row=17; Col=8; Value=3493.875}
    @{Sheet="CUL"; Row=17; Col=9; Value=487.5}
    @{Sheet="CUL"; Row=17; Col=10; Value=6500.25}
    @{Sheet="CUL"; Row=17; Col=11; Value=1462.5}
    @{Sheet="CUL"; Row=17; Col=12; Value=19500.75}
    @{Sheet="CUL"; Row=17; Col=13; Value=-1293.5}
    @{Sheet="CUL"; Row=17; Col=14; Value=-19838.75}
    @{Sheet="CUL"; Row=34; Col=8; Value=393.57144}
    @{Sheet="CUL"; Row=34; Col=10; Value=475.0}
    @{Sheet="CUL"; Row=34; Col=12; Value=1425.0}
    @{Sheet="CUL"; Row=34; Col=14; Value=-1593.0}
    @{Sheet="CUL"; Row=55; Col=8; Value=2583.3333}
    @{Sheet="CUL"; Row=55; Col=10; Value=3125.0}
    @{Sheet="CUL"; Row=55; Col=12; Value=9375.0}
    @{Sheet="CUL"; Row=55; Col=14; Value=-9729.0}
    @{Sheet="GSM"; Row=102; Col=8; Value=7126.5}
    @{Sheet="GSM"; Row=102; Col=9; Value=6170.6665}
    @{Sheet="GSM"; Row=102; Col=10; Value=9994.0}
    @{Sheet="GSM"; Row=102; Col=11; Value=6170.6665}
    @{Sheet="GSM"; Row=102; Col=12; Value=9994.0}
    @{Sheet="GSM"; Row=102; Col=13; Value=-4548.6665}
    @{Sheet="GSM"; Row=102; Col=14; Value=-13238.0}
    @{Sheet="GSM"; Row=132; Col=8; Value=7195.1665}
    @{Sheet="GSM"; Row=132; Col=9; Value=10207.429}
    @{Sheet="GSM"; Row=132; Col=10; Value=2978.0}
    @{Sheet="GSM"; Row=132; Col=11; Value=30622.287}
    @{Sheet="GSM"; Row=132; Col=12; Value=8934.0}
    @{Sheet="GSM"; Row=132; Col=13; Value=-28092.287}
    @{Sheet="GSM"; Row=132; Col=14; Value=-13994.0}
    @{Sheet="LTW"; Row=46; Col=8; Value=1825.0555}
    @{Sheet="LTW"; Row=46; Col=9; Value=1567.1}
    @{Sheet="LTW"; Row=46; Col=10; Value=2147.5}
    @{Sheet="LTW"; Row=46; Col=11; Value=1567.1}
    @{Sheet="LTW"; Row=46; Col=12; Value=2147.5}
    @{Sheet="LTW"; Row=46; Col=13; Value=-1379.1}
    @{Sheet="LTW"; Row=46; Col=14; Value=-2523.5}
    @{Sheet="LTW"; Row=132; Col=8; Value=4323.884}
    @{Sheet="LTW"; Row=132; Col=9; Value=4566.364}
    @{Sheet="LTW"; Row=132; Col=11; Value=13699.092}
    @{Sheet="LTW"; Row=132; Col=13; Value=-11169.092}
    @{Sheet="LTW"; Row=136; Col=8; Value=5680.931}
    @{Sheet="LTW"; Row=136; Col=9; Value=2183.1304}
    @{Sheet="LTW"; Row=136; Col=10; Value=19089.166}
    @{Sheet="LTW"; Row=136; Col=11; Value=6549.3912}
    @{Sheet="LTW"; Row=136; Col=12; Value=57267.49800000001}
    @{Sheet="LTW"; Row=136; Col=13; Value=-3999.3912}
    @{Sheet="LTW"; Row=136; Col=14; Value=-62367.49800000001}
    @{Sheet="WVR"; Row=81; Col=8; Value=720.6429}
    @{Sheet="WVR"; Row=81; Col=9; Value=644.2727}
    @{Sheet="WVR"; Row=81; Col=10; Value=1000.6667}
    @{Sheet="WVR"; Row=81; Col=11; Value=1288.5454}
    @{Sheet="WVR"; Row=81; Col=12; Value=2001.3334}
    @{Sheet="WVR"; Row=81; Col=13; Value=-227.5454}
    @{Sheet="WVR"; Row=81; Col=14; Value=-4123.3334}
    @{Sheet="WVR"; Row=84; Col=8; Value=720.6429}
    @{Sheet="WVR"; Row=84; Col=9; Value=644.2727}
    @{Sheet="WVR"; Row=84; Col=10; Value=1000.6667}
    @{Sheet="WVR"; Row=84; Col=11; Value=6442.727}
    @{Sheet="WVR"; Row=84; Col=12; Value=10006.667}
    @{Sheet="WVR"; Row=84; Col=13; Value=-1138.727}
    @{Sheet="WVR"; Row=84; Col=14; Value=-20614.667}
    @{Sheet="WVR"; Row=132; Col=8; Value=20006136.0}
    @{Sheet="WVR"; Row=132; Col=9; Value=26322922.0}
    @{Sheet="WVR"; Row=132; Col=10; Value=2983.8333}
    @{Sheet="WVR"; Row=132; Col=11; Value=78968766.0}
    @{Sheet="WVR"; Row=132; Col=12; Value=8951.499899999999}
    @{Sheet="WVR"; Row=132; Col=13; Value=-78966236.0}
    @{Sheet="WVR"; Row=132; Col=14; Value=-14011.4999}
    @{Sheet="WVR"; Row=136; Col=8; Value=5299.4443}
    @{Sheet="WVR"; Row=136; Col=9; Value=9647.272}
    @{Sheet="WVR"; Row=136; Col=11; Value=28941.816}
    @{Sheet="WVR"; Row=136; Col=13; Value=-26391.816}
)

foreach ($change in $changes) {
    $ws = $wb.Worksheets.Item($change.Sheet)
    $cell = $ws.Cells.Item($change.Row, $change.Col)
    if ($null -eq $change.Value) {
        $cell.ClearContents()
    } else {
        $cell.Value = $change.Value
    }
}

Write-Host "Applied $($changes.Count) cell updates"
